$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.386.80"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +4.02%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.657.68"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +3.22%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.24%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "567.63"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +6.24%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.21"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.55%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.19%  "

# Row 8
$ws.Range("E8").Value = "  +3.61%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.656.15"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.11%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.85"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.85%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.106"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +6.39%  "

# Row 12
$ws.Range("E12").Value = "  +7.01%  "

# Row 13
$ws.Range("E13").Value = "  +4.34%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.112.69"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.73%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "60.651.40"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.58%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "22.08"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +6.96%  "

# Row 17
$ws.Range("E17").Value = "  +5.83%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.649.53"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.78%  "

# Row 19
$ws.Range("E19").Value = "  +3.74%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "343.05"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.91%  "

# Row 21
$ws.Range("E21").Value = "  +4.51%  "

# Row 22
$ws.Range("E22").Value = "  +4.07%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.99"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.00%  "

# Row 25
$ws.Range("E25").Value = "  +4.69%  "

# Row 26
$ws.Range("E26").Value = "  +5.80%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.995"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.30%  "

# Row 28
$ws.Range("E28").Value = "  +5.75%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0803"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +11.49%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.998"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.06%  "

# Row 31
$ws.Range("E31").Value = "  +4.91%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.16"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.54%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "160.03"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.91%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "19.16"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.07%  "

# Row 35
$ws.Range("E35").Value = "  +6.48%  "

# Row 36
$ws.Range("E36").Value = "  +8.23%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.15"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +6.03%  "

# Row 38
$ws.Range("E38").Value = "  +8.61%  "

# Row 39
$ws.Range("B39").Value = "OKB"
$ws.Range("C39").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "37.62"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.08%  "

# Row 40
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.52"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +7.87%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "299.92"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +6.69%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.64"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.90%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.997"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.27%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0985"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.00%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.605"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.95%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0545"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.48%  "

# Row 47
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "128.82"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +17.63%  "

# Row 48
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "19.46"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +6.45%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "10.71"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.60%  "

# Row 50
$ws.Range("E50").Value = "  +4.50%  "

# Row 51
$ws.Range("B51").Value = "InjectiveProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "18.80"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +6.25%  "
